# DistrictHeating 'smart heating' functionality
# Re-point several grid-connection / actor records from the electricity
# suppliers (sup1/sup2) onto the new Holon actor (hol1), widen the
# "type" column on config_gridConnections, and leave config_actors as the
# active sheet/selection (matching the author's final on-screen state).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "config_gridConnections": I29 sup2(95) -> hol1(99)
# ---------------------------------------------------------------------
$wsGrid = $wb.Worksheets.Item("config_gridConnections")
$wsGrid.Activate()

$wsGrid.Cells.Item(29, 9).Value = "hol1"

# Widen column D ("type") so the new, longer actor name is readable.
$wsGrid.Columns.Item(4).ColumnWidth = 28.15

# Leave the same kind of selection state the author ended up with while
# working on this sheet before moving on to config_actors.
$wsGrid.Range("I32").Select()

# ---------------------------------------------------------------------
# Sheet "config_actors": column E (parent_actor) sup1(94)/sup2(95) -> hol1(99)
# (and a couple of rows revert from hol1 back to sup1)
# ---------------------------------------------------------------------
$wsActors = $wb.Worksheets.Item("config_actors")
$wsActors.Activate()

$wsActors.Cells.Item(27, 5).Value = "hol1"
$wsActors.Cells.Item(28, 5).Value = "sup1"
$wsActors.Cells.Item(34, 5).Value = "hol1"
$wsActors.Cells.Item(35, 5).Value = "hol1"
$wsActors.Cells.Item(36, 5).Value = "hol1"
$wsActors.Cells.Item(45, 5).Value = "hol1"
$wsActors.Cells.Item(52, 5).Value = "hol1"
$wsActors.Cells.Item(53, 5).Value = "hol1"
$wsActors.Cells.Item(54, 5).Value = "hol1"
$wsActors.Cells.Item(63, 5).Value = "hol1"
$wsActors.Cells.Item(70, 5).Value = "hol1"
$wsActors.Cells.Item(71, 5).Value = "hol1"
$wsActors.Cells.Item(72, 5).Value = "hol1"
$wsActors.Cells.Item(81, 5).Value = "hol1"
$wsActors.Cells.Item(88, 5).Value = "hol1"
$wsActors.Cells.Item(89, 5).Value = "hol1"
$wsActors.Cells.Item(90, 5).Value = "hol1"
$wsActors.Cells.Item(99, 5).Value = "hol1"
$wsActors.Cells.Item(106, 5).Value = "hol1"
$wsActors.Cells.Item(107, 5).Value = "hol1"
$wsActors.Cells.Item(110, 5).Value = "hol1"
$wsActors.Cells.Item(117, 5).Value = "hol1"
$wsActors.Cells.Item(118, 5).Value = "hol1"
$wsActors.Cells.Item(119, 5).Value = "hol1"
$wsActors.Cells.Item(128, 5).Value = "hol1"
$wsActors.Cells.Item(135, 5).Value = "hol1"
$wsActors.Cells.Item(136, 5).Value = "hol1"
$wsActors.Cells.Item(137, 5).Value = "hol1"
$wsActors.Cells.Item(146, 5).Value = "hol1"
$wsActors.Cells.Item(153, 5).Value = "hol1"
$wsActors.Cells.Item(154, 5).Value = "hol1"
$wsActors.Cells.Item(155, 5).Value = "hol1"
$wsActors.Cells.Item(164, 5).Value = "hol1"
$wsActors.Cells.Item(171, 5).Value = "hol1"
$wsActors.Cells.Item(172, 5).Value = "hol1"
$wsActors.Cells.Item(173, 5).Value = "hol1"
$wsActors.Cells.Item(182, 5).Value = "hol1"
$wsActors.Cells.Item(189, 5).Value = "hol1"
$wsActors.Cells.Item(190, 5).Value = "hol1"
$wsActors.Cells.Item(193, 5).Value = "hol1"
$wsActors.Cells.Item(194, 5).Value = "sup1"
$wsActors.Cells.Item(200, 5).Value = "hol1"
$wsActors.Cells.Item(201, 5).Value = "hol1"
$wsActors.Cells.Item(202, 5).Value = "hol1"
$wsActors.Cells.Item(211, 5).Value = "hol1"
$wsActors.Cells.Item(218, 5).Value = "hol1"
$wsActors.Cells.Item(219, 5).Value = "hol1"
$wsActors.Cells.Item(220, 5).Value = "hol1"
$wsActors.Cells.Item(229, 5).Value = "hol1"
$wsActors.Cells.Item(236, 5).Value = "hol1"
$wsActors.Cells.Item(237, 5).Value = "hol1"
$wsActors.Cells.Item(238, 5).Value = "hol1"
$wsActors.Cells.Item(247, 5).Value = "hol1"
$wsActors.Cells.Item(254, 5).Value = "hol1"
$wsActors.Cells.Item(255, 5).Value = "hol1"
$wsActors.Cells.Item(256, 5).Value = "hol1"
$wsActors.Cells.Item(265, 5).Value = "hol1"
$wsActors.Cells.Item(272, 5).Value = "hol1"
$wsActors.Cells.Item(273, 5).Value = "hol1"
$wsActors.Cells.Item(274, 5).Value = "hol1"
$wsActors.Cells.Item(277, 5).Value = "sup1"
$wsActors.Cells.Item(283, 5).Value = "hol1"

# Final on-screen state: config_actors active/selected, cell E287 current.
$wsActors.Range("E287").Select()
